# updated ssaform and opt slides
#
# 1) Bump the cached "today" date fields from 2020-11-14 to 2020-11-16
#    (handout master's datetimeFigureOut field, plus the datetime1
#    date-placeholder field cached on the slide master and every
#    slide layout).
# 2) Fix the fully-renamed SSA slide: "4: return k" -> "4: return k2".

$p = $ppt.ActivePresentation

# --- Handout master: M/D/YY style field ("11/14/20" -> "11/16/20") ---
# This placeholder is the auto-updating "datetimeFigureOut" field; PowerPoint
# won't let ordinary text-range edits touch it (same as the real app), so go
# through the Header/Footer dialog's date-and-time setting instead. (The
# DateAndTime.Text getter reflects the dialog's pending value, not the
# placeholder's current text, so check the shape itself first.)
$hm = $p.HandoutMaster
$hmDateShape = $null
for ($i = 1; $i -le $hm.Shapes.Count; $i++) {
    $sh = $hm.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "11/14/20") {
        $hmDateShape = $sh
    }
}
if ($hmDateShape -ne $null) {
    $hm.HeadersFooters.DateAndTime.Text = "11/16/20"
}

# --- Slide master: ISO style field ("2020-11-14" -> "2020-11-16") ---
$sm = $p.SlideMaster
for ($i = 1; $i -le $sm.Shapes.Count; $i++) {
    $sh = $sm.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "2020-11-14") {
            $sh.TextFrame.TextRange.Text = "2020-11-16"
        }
    }
}

# --- Every slide layout carries its own cached copy of the same field ---
for ($j = 1; $j -le $sm.CustomLayouts.Count; $j++) {
    $lay = $sm.CustomLayouts.Item($j)
    for ($i = 1; $i -le $lay.Shapes.Count; $i++) {
        $sh = $lay.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "2020-11-14") {
                $sh.TextFrame.TextRange.Text = "2020-11-16"
            }
        }
    }
}

# --- Slide 14 ("Converting to SSA Form"): return value fix ---
$s14 = $p.Slides.Item(14)
for ($i = 1; $i -le $s14.Shapes.Count; $i++) {
    $sh = $s14.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "4: return k") {
            $sh.TextFrame.TextRange.Text = "4: return k2"
        }
    }
}
